$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data to Sheet1 (A3 = "PC changes again")
$ws.Range("A3").Value = "PC changes again"

# Move / leave the selection on the next empty cell below, as Excel does
# after typing a value into a cell and pressing Enter
$ws.Range("A4").Select()
